$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Ucrania" now ranks ahead of Singapur/Bielorrusia/Catar, so rows 44-47
#     are updated in place (no structural row insertion - Panama in row 48 and
#     beyond stay untouched) ---

# Row 44: Ucrania (new data)
$ws.Range("A44").Value = "Ucrania"
$ws.Range("B44").Value = 5106
$ws.Range("C44").Value = 444
$ws.Range("D44").Value = 275
$ws.Range("E44").Value = 4698
$ws.Range("F44").Value = 45
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 133

# Row 45: Singapur (previous row-44 data)
$ws.Range("A45").Value = "Singapur"
$ws.Range("B45").Value = 5050
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 708
$ws.Range("E45").Value = 4331
$ws.Range("F45").Value = 22
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 11

# Row 46: Bielorrusia (previous row-45 data)
$ws.Range("A46").Value = "Bielorrusia"
$ws.Range("B46").Value = 4779
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 342
$ws.Range("E46").Value = 4395
$ws.Range("F46").Value = 65
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 42

# Row 47: Catar (previous row-46 data)
$ws.Range("A47").Value = "Catar"
$ws.Range("B47").Value = 4663
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 464
$ws.Range("E47").Value = 4192
$ws.Range("F47").Value = 37
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 7

# --- Update Taiwan row (row 109) ---
$ws.Range("B109").Value = 397
$ws.Range("C109").Value = 2
$ws.Range("D109").Value = 178
$ws.Range("E109").Value = 213

# --- Update Georgia row (row 110) ---
$ws.Range("B110").Value = 385
$ws.Range("C110").Value = 15
$ws.Range("D110").Value = 84
$ws.Range("E110").Value = 298
